# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    3 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 12.0302756157461)
    4 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    5 = @(0.06328177979961902, 0.3375848360084654, 0.1529057820181812, 246.9852506941017, 247.5390230919279)
    6 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 6.48142807727062, 12.7228780040422)
    7 = @(3.182878228561681, 87981.0709163148, 2938.103010863317, 6.48142807727062, 90928.83823348394)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Cells.Item($row, 2).Value = $rowVals[0]  # B
    $ws.Cells.Item($row, 3).Value = $rowVals[1]  # C
    $ws.Cells.Item($row, 4).Value = $rowVals[2]  # D
    $ws.Cells.Item($row, 5).Value = $rowVals[3]  # E
    $ws.Cells.Item($row, 7).Value = $rowVals[4]  # G
}
